$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '67.808.28'
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').Value = '3.805.93'
$ws.Range('E3').Value = '  +0.21%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.15%  '
Set-TextValue 'D5' '596.93'
$ws.Range('E5').Value = '  +0.13%  '
Set-TextValue 'D6' '167.05'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').Value = '3.803.19'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  +0.10%  '
Set-TextValue 'D11' '6.30'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('E12').Value = '  -0.70%  '
Set-TextValue 'D13' '0.0000254'
$ws.Range('E13').Value = '  -2.65%  '
Set-TextValue 'D14' '36.11'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '4.443.48'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').Value = '3.840.86'
$ws.Range('E16').Value = '  +1.31%  '
Set-TextValue 'D17' '18.64'
$ws.Range('E17').Value = '  +4.56%  '
$ws.Range('D18').Value = '67.791.13'
$ws.Range('E18').Value = '  -0.89%  '
Set-TextValue 'D19' '7.12'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('E20').Value = '  +0.22%  '
Set-TextValue 'D21' '461.50'
$ws.Range('E21').Value = '  -0.61%  '
Set-TextValue 'D22' '9.88'
$ws.Range('E22').Value = '  -7.80%  '
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  -0.83%  '
Set-TextValue 'D25' '83.52'
$ws.Range('E25').Value = '  -0.51%  '
Set-TextValue 'D26' '12.08'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('E27').Value = '  -2.87%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D28' '10.02'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').Value = '3.950.48'
$ws.Range('E30').Value = '  -0.09%  '
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('E32').Value = '  +4.12%  '
Set-TextValue 'D33' '7.25'
$ws.Range('E33').Value = '  -0.55%  '
Set-TextValue 'D34' '29.72'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('E35').Value = '  +0.02%  '
Set-TextValue 'D36' '9.08'
$ws.Range('E36').Value = '  -0.72%  '
Set-TextValue 'D37' '0.100'
$ws.Range('E37').Value = '  -0.49%  '
Set-TextValue 'D38' '3.37'
$ws.Range('E38').Value = '  -3.69%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  -0.29%  '
Set-TextValue 'D41' '5.79'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +0.06%  '
Set-TextValue 'D44' '48.12'
$ws.Range('E44').Value = '  +2.78%  '
Set-TextValue 'D45' '43.93'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('E46').Value = '  -1.18%  '
Set-TextValue 'D47' '150.85'
$ws.Range('E47').Value = '  +2.65%  '
Set-TextValue 'D48' '8.33'
$ws.Range('E48').Value = '  -1.06%  '
Set-TextValue 'D49' '392.45'
$ws.Range('E49').Value = '  +0.54%  '
$ws.Range('E50').Value = '  -4.68%  '
Set-TextValue 'D51' '26.39'
$ws.Range('E51').Value = '  +5.01%  '
